$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-24 Thursday" "2024-10-25 Friday"

Replace-Text "862÷4=" "205÷7="
Replace-Text "510÷6=" "746÷2="
Replace-Text "868÷3=" "112÷7="
Replace-Text "118÷8=" "548÷3="
Replace-Text "318÷6=" "445÷5="
Replace-Text "700÷3=" "182÷8="
Replace-Text "789÷7=" "837÷5="
Replace-Text "944÷7=" "945÷3="
Replace-Text "889÷7=" "466÷2="
Replace-Text "957÷8=" "667÷6="
Replace-Text "320÷8=" "908÷7="
Replace-Text "609÷7=" "528÷4="
Replace-Text "258÷4=" "864÷3="
Replace-Text "131÷8=" "620÷6="
Replace-Text "878÷4=" "499÷8="
Replace-Text "313÷5=" "330÷8="
Replace-Text "701÷6=" "436÷3="
Replace-Text "298÷2=" "705÷2="
Replace-Text "106÷2=" "268÷9="
Replace-Text "651÷7=" "224÷6="
Replace-Text "999÷2=" "952÷4="
Replace-Text "354÷3=" "847÷4="
Replace-Text "271÷5=" "303÷3="
Replace-Text "290÷7=" "806÷7="
Replace-Text "511÷9=" "319÷9="
